$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 632.03705
$ws.Range("I53").Value = 592.6667
$ws.Range("J53").Value = 710.7778
$ws.Range("K53").Value = 592.6667
$ws.Range("L53").Value = 710.7778
$ws.Range("M53").Value = 44.33330000000001
$ws.Range("N53").Value = -1984.7778
$ws.Range("H101").Value = 4017.1667
$ws.Range("I101").Value = 1651.4
$ws.Range("J101").Value = 6974.375
$ws.Range("K101").Value = 4954.200000000001
$ws.Range("L101").Value = 20923.125
$ws.Range("M101").Value = -3332.200000000001
$ws.Range("N101").Value = -24167.125
$ws.Range("H106").Value = 3360.5557
$ws.Range("I106").Value = 3360.5557
$ws.Range("K106").Value = 3360.5557
$ws.Range("M106").Value = -2729.5557
$ws.Range("H138").Value = 4718.724
$ws.Range("J138").Value = 4489.706
$ws.Range("L138").Value = 13469.118
$ws.Range("N138").Value = -23749.118

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3119.85
$ws.Range("I61").Value = 3020.8948
$ws.Range("K61").Value = 3020.8948
$ws.Range("M61").Value = -2808.8948
$ws.Range("H74").Value = 3148.25
$ws.Range("J74").Value = 6994
$ws.Range("L74").Value = 6994
$ws.Range("N74").Value = -8742
$ws.Range("H77").Value = 3148.25
$ws.Range("J77").Value = 6994
$ws.Range("L77").Value = 34970
$ws.Range("N77").Value = -43706
$ws.Range("H132").Value = 2382.4614
$ws.Range("I132").Value = 2401.2727
$ws.Range("K132").Value = 7203.8181
$ws.Range("M132").Value = -4673.8181
$ws.Range("H136").Value = 3119.85
$ws.Range("I136").Value = 3020.8948
$ws.Range("K136").Value = 9062.6844
$ws.Range("M136").Value = -6512.6844

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1120.7142
$ws.Range("I94").Value = 1120.7142
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1120.7142
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -669.7141999999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 49816.5
$ws.Range("J9").Value = 49816.5
$ws.Range("L9").Value = 49816.5
$ws.Range("N9").Value = -50152.5
$ws.Range("H12").Value = 2633.2
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 5833
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 5833
$ws.Range("M12").Value = -330
$ws.Range("N12").Value = -6173
$ws.Range("H31").Value = 4821.9443
$ws.Range("I31").Value = 3228
$ws.Range("K31").Value = 3228
$ws.Range("M31").Value = -2933
$ws.Range("H34").Value = 4821.9443
$ws.Range("I34").Value = 3228
$ws.Range("K34").Value = 3228
$ws.Range("M34").Value = -3026
$ws.Range("H99").Value = 5842.2856
$ws.Range("I99").Value = 6174.5
$ws.Range("J99").Value = 5399.3335
$ws.Range("K99").Value = 6174.5
$ws.Range("L99").Value = 5399.3335
$ws.Range("M99").Value = -4676.5
$ws.Range("N99").Value = -8395.333500000001
$ws.Range("H104").Value = 30285
$ws.Range("J104").Value = 30285
$ws.Range("L104").Value = 30285
$ws.Range("N104").Value = -35527
$ws.Range("H126").Value = 5842.2856
$ws.Range("I126").Value = 6174.5
$ws.Range("J126").Value = 5399.3335
$ws.Range("K126").Value = 18523.5
$ws.Range("L126").Value = 16198.0005
$ws.Range("M126").Value = -16053.5
$ws.Range("N126").Value = -21138.0005
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280
$ws.Range("H141").Value = 160108.31
$ws.Range("J141").Value = 173444.8
$ws.Range("L141").Value = 173444.8
$ws.Range("N141").Value = -183804.8

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 999.6667
$ws.Range("I3").Value = 999.6667
$ws.Range("K3").Value = 2999.0001
$ws.Range("M3").Value = -2887.0001
$ws.Range("H11").Value = 11456064
$ws.Range("I11").Value = 18001258
$ws.Range("J11").Value = 1973.5
$ws.Range("K11").Value = 54003774
$ws.Range("L11").Value = 5920.5
$ws.Range("M11").Value = -54003634
$ws.Range("N11").Value = -6200.5
$ws.Range("H31").Value = 830
$ws.Range("J31").Value = 550
$ws.Range("L31").Value = 1650
$ws.Range("N31").Value = -2226
$ws.Range("H132").Value = 1594.7028
$ws.Range("J132").Value = 1694.2354
$ws.Range("L132").Value = 15248.1186
$ws.Range("N132").Value = -20308.1186

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 51.52381
$ws.Range("I2").Value = 49.944443
$ws.Range("J2").Value = 61
$ws.Range("K2").Value = 49.944443
$ws.Range("L2").Value = 61
$ws.Range("M2").Value = 63.055557
$ws.Range("N2").Value = -287
$ws.Range("H62").Value = 41099.75
$ws.Range("I62").Value = 41099.75
$ws.Range("K62").Value = 41099.75
$ws.Range("M62").Value = -40413.75
$ws.Range("H65").Value = 41099.75
$ws.Range("I65").Value = 41099.75
$ws.Range("K65").Value = 123299.25
$ws.Range("M65").Value = -119867.25
$ws.Range("H70").Value = 6820.6665
$ws.Range("I70").Value = 4123
$ws.Range("J70").Value = 9788.1
$ws.Range("K70").Value = 4123
$ws.Range("L70").Value = 9788.1
$ws.Range("M70").Value = -3853
$ws.Range("N70").Value = -10328.1
$ws.Range("H73").Value = 6820.6665
$ws.Range("I73").Value = 4123
$ws.Range("J73").Value = 9788.1
$ws.Range("K73").Value = 4123
$ws.Range("L73").Value = 9788.1
$ws.Range("M73").Value = -3187
$ws.Range("N73").Value = -11660.1
$ws.Range("H113").Value = 801447
$ws.Range("I113").Value = 801447
$ws.Range("K113").Value = 801447
$ws.Range("M113").Value = -799277
$ws.Range("H126").Value = 4654.077
$ws.Range("I126").Value = 3844.889
$ws.Range("K126").Value = 11534.667
$ws.Range("M126").Value = -9064.667000000001
$ws.Range("H137").Value = 62671.875
$ws.Range("I137").Value = 62671.875
$ws.Range("K137").Value = 62671.875
$ws.Range("M137").Value = -57571.875

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1438
$ws.Range("J22").Value = 1249.5834
$ws.Range("L22").Value = 1249.5834
$ws.Range("N22").Value = -1839.5834
$ws.Range("H27").Value = 1438
$ws.Range("J27").Value = 1249.5834
$ws.Range("L27").Value = 1249.5834
$ws.Range("N27").Value = -1463.5834
$ws.Range("H132").Value = 13395.342
$ws.Range("I132").Value = 13283.345
$ws.Range("K132").Value = 39850.035
$ws.Range("M132").Value = -37320.035
$ws.Range("H136").Value = 7819.364
$ws.Range("J136").Value = 6333.3335
$ws.Range("L136").Value = 19000.0005
$ws.Range("N136").Value = -24100.0005

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 10494.667
$ws.Range("J3").Value = 742.5
$ws.Range("L3").Value = 742.5
$ws.Range("N3").Value = -970.5
$ws.Range("H6").Value = 3350
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 5700
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 5700
$ws.Range("M6").Value = -885
$ws.Range("N6").Value = -5930
$ws.Range("H15").Value = 7266.6665
$ws.Range("J15").Value = 7266.6665
$ws.Range("L15").Value = 7266.6665
$ws.Range("N15").Value = -7842.6665
$ws.Range("H41").Value = 18824.154
$ws.Range("I41").Value = 18514.5
$ws.Range("J41").Value = 18961.777
$ws.Range("K41").Value = 18514.5
$ws.Range("L41").Value = 18961.777
$ws.Range("M41").Value = -18124.5
$ws.Range("N41").Value = -19741.777
$ws.Range("H132").Value = 4743.1665
$ws.Range("I132").Value = 4357.6113
$ws.Range("J132").Value = 5899.8335
$ws.Range("K132").Value = 13072.8339
$ws.Range("L132").Value = 17699.5005
$ws.Range("M132").Value = -10542.8339
$ws.Range("N132").Value = -22759.5005
$ws.Range("H135").Value = 76663.664
$ws.Range("J135").Value = 76663.664
$ws.Range("L135").Value = 76663.664
$ws.Range("N135").Value = -86803.664
$ws.Range("H136").Value = 3815.535
$ws.Range("I136").Value = 2457.8823
$ws.Range("J136").Value = 8944.444
$ws.Range("K136").Value = 7373.646900000001
$ws.Range("L136").Value = 26833.332
$ws.Range("M136").Value = -4823.646900000001
$ws.Range("N136").Value = -31933.332
